# Auto-generated Excel COM-interop edit script
# Updates the 3 worksheets (LP1912, LP1912-215, 6203-6173) with the refreshed
# scraped schedule data: new 'Ultima actualizacion' / 'Total filas' headers,
# and the full data table (rows 6+) reflecting the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:50'
$ws.Cells.Item(3, 1).Value = 'Total filas: 55'

$ws.Cells.Item(29, 1).Value = '05:54:50'
$ws.Cells.Item(29, 2).Value = '05:54'
$ws.Cells.Item(29, 3).Value = '10_OLMOS'
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 'LP1912'

$ws.Cells.Item(30, 1).Value = '05:54:50'
$ws.Cells.Item(30, 2).Value = '06:04'
$ws.Cells.Item(30, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(30, 4).Value = 10
$ws.Cells.Item(30, 5).Value = 'LP1912'

$ws.Cells.Item(32, 1).Value = '05:54:50'
$ws.Cells.Item(32, 2).Value = '06:11'
$ws.Cells.Item(32, 3).Value = '215A_EL PATO'
$ws.Cells.Item(32, 4).Value = 17
$ws.Cells.Item(32, 5).Value = 'LP1912'

$ws.Cells.Item(34, 1).Value = '05:54:50'
$ws.Cells.Item(34, 2).Value = '06:14'
$ws.Cells.Item(34, 3).Value = '225_HARAS DEL SUR'
$ws.Cells.Item(34, 4).Value = 20
$ws.Cells.Item(34, 5).Value = 'LP1912'

$ws.Cells.Item(35, 1).Value = '05:54:50'
$ws.Cells.Item(35, 2).Value = '06:21'
$ws.Cells.Item(35, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(35, 4).Value = 27
$ws.Cells.Item(35, 5).Value = 'LP1912'

$ws.Cells.Item(36, 1).Value = '05:54:50'
$ws.Cells.Item(36, 2).Value = '06:27'
$ws.Cells.Item(36, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(36, 4).Value = 33
$ws.Cells.Item(36, 5).Value = 'LP1912'

$ws.Cells.Item(37, 1).Value = '05:54:50'
$ws.Cells.Item(37, 2).Value = '06:29'
$ws.Cells.Item(37, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(37, 4).Value = 35
$ws.Cells.Item(37, 5).Value = 'LP1912'

$ws.Cells.Item(40, 1).Value = '05:54:50'
$ws.Cells.Item(40, 2).Value = '06:44'
$ws.Cells.Item(40, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(40, 4).Value = 50
$ws.Cells.Item(40, 5).Value = 'LP1912'

$ws.Cells.Item(41, 1).Value = '05:54:50'
$ws.Cells.Item(41, 2).Value = '06:46'
$ws.Cells.Item(41, 3).Value = '215C_EL PATO'
$ws.Cells.Item(41, 4).Value = 52
$ws.Cells.Item(41, 5).Value = 'LP1912'

$ws.Cells.Item(43, 1).Value = '05:54:50'
$ws.Cells.Item(43, 2).Value = '06:59'
$ws.Cells.Item(43, 3).Value = '14_ABASTO'
$ws.Cells.Item(43, 4).Value = 65
$ws.Cells.Item(43, 5).Value = 'LP1912'

$ws.Cells.Item(44, 1).Value = '05:23:05'
$ws.Cells.Item(44, 2).Value = '07:00'
$ws.Cells.Item(44, 3).Value = '14_ABASTO'
$ws.Cells.Item(44, 4).Value = 97
$ws.Cells.Item(44, 5).Value = 'LP1912'

$ws.Cells.Item(45, 1).Value = '05:54:50'
$ws.Cells.Item(45, 2).Value = '07:04'
$ws.Cells.Item(45, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(45, 4).Value = 70
$ws.Cells.Item(45, 5).Value = 'LP1912'

$ws.Cells.Item(46, 1).Value = '05:54:50'
$ws.Cells.Item(46, 2).Value = '07:05'
$ws.Cells.Item(46, 3).Value = '15_ABASTO'
$ws.Cells.Item(46, 4).Value = 71
$ws.Cells.Item(46, 5).Value = 'LP1912'

$ws.Cells.Item(47, 1).Value = '05:54:50'
$ws.Cells.Item(47, 2).Value = '07:06'
$ws.Cells.Item(47, 3).Value = '225_GOMEZ'
$ws.Cells.Item(47, 4).Value = 72
$ws.Cells.Item(47, 5).Value = 'LP1912'

$ws.Cells.Item(48, 1).Value = '05:23:05'
$ws.Cells.Item(48, 2).Value = '07:07'
$ws.Cells.Item(48, 3).Value = '225_GOMEZ'
$ws.Cells.Item(48, 4).Value = 104
$ws.Cells.Item(48, 5).Value = 'LP1912'

$ws.Cells.Item(49, 1).Value = '05:54:50'
$ws.Cells.Item(49, 2).Value = '07:11'
$ws.Cells.Item(49, 3).Value = '215A_EL PATO'
$ws.Cells.Item(49, 4).Value = 77
$ws.Cells.Item(49, 5).Value = 'LP1912'

$ws.Cells.Item(50, 1).Value = '05:23:05'
$ws.Cells.Item(50, 2).Value = '07:12'
$ws.Cells.Item(50, 3).Value = '215A_EL PATO'
$ws.Cells.Item(50, 4).Value = 109
$ws.Cells.Item(50, 5).Value = 'LP1912'

$ws.Cells.Item(51, 1).Value = '05:54:50'
$ws.Cells.Item(51, 2).Value = '07:15'
$ws.Cells.Item(51, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(51, 4).Value = 81
$ws.Cells.Item(51, 5).Value = 'LP1912'

$ws.Cells.Item(52, 1).Value = '05:23:05'
$ws.Cells.Item(52, 2).Value = '07:16'
$ws.Cells.Item(52, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(52, 4).Value = 113
$ws.Cells.Item(52, 5).Value = 'LP1912'

$ws.Cells.Item(53, 1).Value = '05:54:50'
$ws.Cells.Item(53, 2).Value = '07:21'
$ws.Cells.Item(53, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(53, 4).Value = 87
$ws.Cells.Item(53, 5).Value = 'LP1912'

$ws.Cells.Item(54, 1).Value = '05:54:50'
$ws.Cells.Item(54, 2).Value = '07:23'
$ws.Cells.Item(54, 3).Value = '10_OLMOS'
$ws.Cells.Item(54, 4).Value = 89
$ws.Cells.Item(54, 5).Value = 'LP1912'

$ws.Cells.Item(55, 1).Value = '05:54:50'
$ws.Cells.Item(55, 2).Value = '07:31'
$ws.Cells.Item(55, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(55, 4).Value = 97
$ws.Cells.Item(55, 5).Value = 'LP1912'

$ws.Cells.Item(56, 1).Value = '05:54:50'
$ws.Cells.Item(56, 2).Value = '07:32'
$ws.Cells.Item(56, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(56, 4).Value = 98
$ws.Cells.Item(56, 5).Value = 'LP1912'

$ws.Cells.Item(57, 1).Value = '05:54:50'
$ws.Cells.Item(57, 2).Value = '07:36'
$ws.Cells.Item(57, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(57, 4).Value = 102
$ws.Cells.Item(57, 5).Value = 'LP1912'

$ws.Cells.Item(58, 1).Value = '05:54:50'
$ws.Cells.Item(58, 2).Value = '07:46'
$ws.Cells.Item(58, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(58, 4).Value = 112
$ws.Cells.Item(58, 5).Value = 'LP1912'

$ws.Cells.Item(59, 1).Value = '05:54:50'
$ws.Cells.Item(59, 2).Value = '07:47'
$ws.Cells.Item(59, 3).Value = '14_ABASTO'
$ws.Cells.Item(59, 4).Value = 113
$ws.Cells.Item(59, 5).Value = 'LP1912'

$ws.Cells.Item(60, 1).Value = '05:54:50'
$ws.Cells.Item(60, 2).Value = '07:51'
$ws.Cells.Item(60, 3).Value = '215D_EL PATO'
$ws.Cells.Item(60, 4).Value = 117
$ws.Cells.Item(60, 5).Value = 'LP1912'


# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:50'
$ws.Cells.Item(3, 1).Value = 'Total filas: 16'

$ws.Cells.Item(15, 1).Value = '05:54:50'
$ws.Cells.Item(15, 2).Value = '06:11'
$ws.Cells.Item(15, 3).Value = '215A_EL PATO'
$ws.Cells.Item(15, 4).Value = 17
$ws.Cells.Item(15, 5).Value = 'LP1912'

$ws.Cells.Item(17, 1).Value = '05:54:50'
$ws.Cells.Item(17, 2).Value = '06:46'
$ws.Cells.Item(17, 3).Value = '215C_EL PATO'
$ws.Cells.Item(17, 4).Value = 52
$ws.Cells.Item(17, 5).Value = 'LP1912'

$ws.Cells.Item(19, 1).Value = '05:54:50'
$ws.Cells.Item(19, 2).Value = '07:11'
$ws.Cells.Item(19, 3).Value = '215A_EL PATO'
$ws.Cells.Item(19, 4).Value = 77
$ws.Cells.Item(19, 5).Value = 'LP1912'

$ws.Cells.Item(20, 1).Value = '05:23:05'
$ws.Cells.Item(20, 2).Value = '07:12'
$ws.Cells.Item(20, 3).Value = '215A_EL PATO'
$ws.Cells.Item(20, 4).Value = 109
$ws.Cells.Item(20, 5).Value = 'LP1912'

$ws.Cells.Item(21, 1).Value = '05:54:50'
$ws.Cells.Item(21, 2).Value = '07:51'
$ws.Cells.Item(21, 3).Value = '215D_EL PATO'
$ws.Cells.Item(21, 4).Value = 117
$ws.Cells.Item(21, 5).Value = 'LP1912'


# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 05:54:50'
$ws.Cells.Item(3, 1).Value = 'Total filas: 9'

$ws.Cells.Item(8, 1).Value = '05:54:50'
$ws.Cells.Item(8, 2).Value = '06:08'
$ws.Cells.Item(8, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(8, 4).Value = 14
$ws.Cells.Item(8, 5).Value = 'L6173'

$ws.Cells.Item(9, 1).Value = '05:23:05'
$ws.Cells.Item(9, 2).Value = '06:09'
$ws.Cells.Item(9, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(9, 4).Value = 46
$ws.Cells.Item(9, 5).Value = 'L6173'

$ws.Cells.Item(10, 1).Value = '05:54:50'
$ws.Cells.Item(10, 2).Value = '06:32'
$ws.Cells.Item(10, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(10, 4).Value = 38
$ws.Cells.Item(10, 5).Value = 'L6203'

$ws.Cells.Item(11, 1).Value = '05:23:05'
$ws.Cells.Item(11, 2).Value = '06:33'
$ws.Cells.Item(11, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(11, 4).Value = 70
$ws.Cells.Item(11, 5).Value = 'L6203'

$ws.Cells.Item(12, 1).Value = '05:54:50'
$ws.Cells.Item(12, 2).Value = '06:59'
$ws.Cells.Item(12, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(12, 4).Value = 65
$ws.Cells.Item(12, 5).Value = 'L6173'

$ws.Cells.Item(13, 1).Value = '05:23:05'
$ws.Cells.Item(13, 2).Value = '07:00'
$ws.Cells.Item(13, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(13, 4).Value = 97
$ws.Cells.Item(13, 5).Value = 'L6173'

$ws.Cells.Item(14, 1).Value = '05:54:50'
$ws.Cells.Item(14, 2).Value = '07:35'
$ws.Cells.Item(14, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(14, 4).Value = 101
$ws.Cells.Item(14, 5).Value = 'L6173'

